$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the column definitions: column A should only be min=1 max=1 (it was
#    incorrectly overlapping column B's width definition at min=1 max=2).
#    Touching column B's width forces the engine to split the old A:B
#    column-format group, leaving column A with its own clean, exact
#    min=1/max=1 definition.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# ---------------------------------------------------------------------------
# 2. Insert two new blank rows at position 13, pushing the old rows 13-24
#    down to 15-26 (row heights travel with the shifted rows automatically).
# ---------------------------------------------------------------------------
$ws.Range("A13:A14").EntireRow.Insert()

# The insert leaves behind empty formatted cells in column A for the two new
# rows (carried over from row 13's old formatting); the target layout has no
# A13/A14 cells at all, so clear them out completely.
$ws.Range("A13:A14").Clear()

# ---------------------------------------------------------------------------
# Helper text blocks (kept as variables for readability).
# ---------------------------------------------------------------------------
$objetivosTxt = "Apresentar as técnicas de caracterização e processamento de polímeros de forma que aluno seja capaz de compreender a importância de cada uma e relacioná-las com o papel desempenhado pelos polímeros na Engenharia de Materiais."
$programaResumidoTxt = "Caracterização de Polímeros: Identificação de polímeros; Determinação da massa molar média de polímeros; Determinação das propriedades físicas; análise térmica de polímeros. Processamento de materiais poliméricos"
$programaTxt = "Testes simples: queima e densidade relativa; FTIR de polímeros; Princípios de Ressonância Magnética Nuclear aplicada a polímeros; Viscosimetria e Reologia; Cromatografia por Exclusão de Tamanho (SEC/GPC); Termogravimetria e Calorimetria Diferencial Exploratória (DSC) de polímeros; Análise Dinâmico – Mecânica (DMA). Índice de Fluidez. Moagem e moldagem de polímeros."
$bibliografiaTxt = "HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill Inc, 1992S. V. CANEVAROLO Jr. Técnicas de Caracterização de Polímeros. São Paulo: Editora Artliber, 2005.MANRICH, S. Processamento de Termoplásticos. Editora Artliber, 2005. NAVARRO, R.F. Fundamentos de Reologia de Polímeros. Editora da Universidade de Caxias do Sul, 1997. MANO, E. B.; MENDES, L. C. Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000. TURI, E. A. Thermal Characterization of Polymeric Materials. New York: Academic Press, 1981.NAVARRO, R.F. Fundamentos de Reologia de Polímeros. Editora da Universidade de Caxias do Sul, 1997.MANO, E. B.; MENDES, L. C. Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000."

function Set-BC($row, $text) {
    $b = $ws.Cells.Item($row, 2)
    $b.Value = $text
    $b.Font.Bold = $false
    $b.WrapText = $true
    $b.VerticalAlignment = -4160

    $c = $ws.Cells.Item($row, 3)
    $c.Value = $text
    $c.Font.Bold = $false
    $c.Font.Color = 255
    $c.WrapText = $true
    $c.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 3. Row 10 ("Objetivos:") - replace the wrongly duplicated professor name
#    with the real objectives paragraph.
# ---------------------------------------------------------------------------
Set-BC 10 $objetivosTxt

# ---------------------------------------------------------------------------
# 4. The two freshly inserted rows (13 & 14) hold the professors' names that
#    used to (incorrectly) live elsewhere.
# ---------------------------------------------------------------------------
Set-BC 13 "5840897 - Clodoaldo Saron"
Set-BC 14 "1033242 - Fábio Herbst Florenzano"

# ---------------------------------------------------------------------------
# 5. Row 15 ("Programa resumido:") gets its real summary text (previously
#    held a stray date value).
# ---------------------------------------------------------------------------
Set-BC 15 $programaResumidoTxt

# ---------------------------------------------------------------------------
# 6. Row 17 ("Programa:") gets its real syllabus text (previously held a
#    stray professor name).
# ---------------------------------------------------------------------------
Set-BC 17 $programaTxt

# ---------------------------------------------------------------------------
# 7. Row 23 ("Bibliografia:") gets the real bibliography text (previously
#    held the "recuperação" text that has now moved up to row 22).
# ---------------------------------------------------------------------------
Set-BC 23 $bibliografiaTxt

Write-Host "done"
